$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# Proofread corrections to the description text in the sample data.
$ws.Range("A7").Value = "This row is invalid and will be skipped. The other rows should yield a total of 10 manifestations."
$ws.Range("C6").Value = 'Should yield 4 (2*2) manifestations AFTER first being unable to sort due to the "1" not being in quotes.'

# Move the active selection back to A2 on the Entities sheet.
[void]$ws.Range("A2").Select()
